$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cliente")

# Row 13 corresponds to "VENDA 09" sale entry.
$ws.Range("B13").Value = "SILVIO TADEU"
$ws.Range("C13").Value = "350b8ddaabcd000fb6470ca8c0e11441"
$ws.Range("D13").Value = Get-Date -Year 2022 -Month 10 -Day 3 -Hour 0 -Minute 0 -Second 0
$ws.Range("G13").Value = "VENDA 09 (03/10)"
